$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.72262366666667
$ws.Range("H2").Value = 89.167871
$ws.Range("I2").Value = 0.2181294538094716
$ws.Range("J2").Value = 0.2181294538094716
$ws.Range("M2").Value = 0.01339666666666667
$ws.Range("N2").Value = 0.04019
$ws.Range("O2").Value = 0.08393217762128816
$ws.Range("P2").Value = 0.08393217762128814
$ws.Range("Q2").Value = 0.3981840817211111
$ws.Range("R2").Value = 3.58365673549
$ws.Range("S2").Value = 0.01830808006157114
$ws.Range("T2").Value = 0.01830808006157114

$ws.Range("G3").Value = 29.72262366666667
$ws.Range("H3").Value = 89.167871
$ws.Range("I3").Value = 0.2181294538094716
$ws.Range("J3").Value = 0.2181294538094716
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.1260863333333333
$ws.Range("N3").Value = 0.378259
$ws.Range("O3").Value = 0.7899502755623498
$ws.Range("P3").Value = 0.7899502755623498
$ws.Range("Q3").Value = 3.747616635176556
$ws.Range("R3").Value = 33.728549716589
$ws.Range("S3").Value = 0.172311422145057
$ws.Range("T3").Value = 0.172311422145057

$ws.Range("G4").Value = 29.72262366666667
$ws.Range("H4").Value = 89.167871
$ws.Range("I4").Value = 0.2181294538094716
$ws.Range("J4").Value = 0.2181294538094716
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02013
$ws.Range("N4").Value = 0.06039
$ws.Range("O4").Value = 0.1261175468163621
$ws.Range("P4").Value = 0.1261175468163621
$ws.Range("Q4").Value = 0.5983164144099999
$ws.Range("R4").Value = 5.384847729690001
$ws.Range("S4").Value = 0.02750995160284352
$ws.Range("T4").Value = 0.02750995160284352

$ws.Range("I5").Value = 0.1680996954419849
$ws.Range("J5").Value = 0.1680996954419849
$ws.Range("M5").Value = 0.01339666666666667
$ws.Range("N5").Value = 0.04019
$ws.Range("O5").Value = 0.08393217762128816
$ws.Range("P5").Value = 0.08393217762128814
$ws.Range("Q5").Value = 0.3068573349366667
$ws.Range("R5").Value = 2.76171601443
$ws.Range("S5").Value = 0.01410897349592112
$ws.Range("T5").Value = 0.01410897349592112

$ws.Range("I6").Value = 0.1680996954419849
$ws.Range("J6").Value = 0.1680996954419849
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1260863333333333
$ws.Range("N6").Value = 0.378259
$ws.Range("O6").Value = 0.7899502755623498
$ws.Range("P6").Value = 0.7899502755623498
$ws.Range("Q6").Value = 2.888070382080334
$ws.Range("R6").Value = 25.992633438723
$ws.Range("S6").Value = 0.132790400736343
$ws.Range("T6").Value = 0.132790400736343

$ws.Range("I7").Value = 0.1680996954419849
$ws.Range("J7").Value = 0.1680996954419849
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02013
$ws.Range("N7").Value = 0.06039
$ws.Range("O7").Value = 0.1261175468163621
$ws.Range("P7").Value = 0.1261175468163621
$ws.Range("Q7").Value = 0.46108769487
$ws.Range("R7").Value = 4.14978925383
$ws.Range("S7").Value = 0.02120032120972073
$ws.Range("T7").Value = 0.02120032120972073

$ws.Range("G8").Value = 44.20756266666667
$ws.Range("H8").Value = 132.622688
$ws.Range("I8").Value = 0.3244320423012451
$ws.Range("J8").Value = 0.3244320423012451
$ws.Range("M8").Value = 0.01339666666666667
$ws.Range("N8").Value = 0.04019
$ws.Range("O8").Value = 0.08393217762128816
$ws.Range("P8").Value = 0.08393217762128814
$ws.Range("Q8").Value = 0.5922339811911111
$ws.Range("R8").Value = 5.33010583072
$ws.Range("S8").Value = 0.02723028780046538
$ws.Range("T8").Value = 0.02723028780046537

$ws.Range("G9").Value = 44.20756266666667
$ws.Range("H9").Value = 132.622688
$ws.Range("I9").Value = 0.3244320423012451
$ws.Range("J9").Value = 0.3244320423012451
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.1260863333333333
$ws.Range("N9").Value = 0.378259
$ws.Range("O9").Value = 0.7899502755623498
$ws.Range("P9").Value = 0.7899502755623498
$ws.Range("Q9").Value = 5.573969482243555
$ws.Range("R9").Value = 50.165725340192
$ws.Range("S9").Value = 0.2562851812171245
$ws.Range("T9").Value = 0.2562851812171245

$ws.Range("G10").Value = 44.20756266666667
$ws.Range("H10").Value = 132.622688
$ws.Range("I10").Value = 0.3244320423012451
$ws.Range("J10").Value = 0.3244320423012451
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02013
$ws.Range("N10").Value = 0.06039
$ws.Range("O10").Value = 0.1261175468163621
$ws.Range("P10").Value = 0.1261175468163621
$ws.Range("Q10").Value = 0.88989823648
$ws.Range("R10").Value = 8.009084128320001
$ws.Range("S10").Value = 0.04091657328365524
$ws.Range("T10").Value = 0.04091657328365524

$ws.Range("G11").Value = 12.94423466666667
$ws.Range("H11").Value = 38.832704
$ws.Range("I11").Value = 0.09499561241587662
$ws.Range("J11").Value = 0.09499561241587662
$ws.Range("M11").Value = 0.01339666666666667
$ws.Range("N11").Value = 0.04019
$ws.Range("O11").Value = 0.08393217762128816
$ws.Range("P11").Value = 0.08393217762128814
$ws.Range("Q11").Value = 0.1734095970844444
$ws.Range("R11").Value = 1.56068637376
$ws.Range("S11").Value = 0.007973188614532403
$ws.Range("T11").Value = 0.007973188614532402

$ws.Range("G12").Value = 12.94423466666667
$ws.Range("H12").Value = 38.832704
$ws.Range("I12").Value = 0.09499561241587662
$ws.Range("J12").Value = 0.09499561241587662
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.1260863333333333
$ws.Range("N12").Value = 0.378259
$ws.Range("O12").Value = 0.7899502755623498
$ws.Range("P12").Value = 0.7899502755623498
$ws.Range("Q12").Value = 1.632091086926222
$ws.Range("R12").Value = 14.688819782336
$ws.Range("S12").Value = 0.07504181020513591
$ws.Range("T12").Value = 0.07504181020513591

$ws.Range("G13").Value = 12.94423466666667
$ws.Range("H13").Value = 38.832704
$ws.Range("I13").Value = 0.09499561241587662
$ws.Range("J13").Value = 0.09499561241587662
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.02013
$ws.Range("N13").Value = 0.06039
$ws.Range("O13").Value = 0.1261175468163621
$ws.Range("P13").Value = 0.1261175468163621
$ws.Range("Q13").Value = 0.26056744384
$ws.Range("R13").Value = 2.34510699456
$ws.Range("S13").Value = 0.0119806135962083
$ws.Range("T13").Value = 0.0119806135962083

$ws.Range("G14").Value = 26.481475
$ws.Range("H14").Value = 79.444425
$ws.Range("I14").Value = 0.1943431960314218
$ws.Range("J14").Value = 0.1943431960314218
$ws.Range("M14").Value = 0.01339666666666667
$ws.Range("N14").Value = 0.04019
$ws.Range("O14").Value = 0.08393217762128816
$ws.Range("P14").Value = 0.08393217762128814
$ws.Range("Q14").Value = 0.3547634934166666
$ws.Range("R14").Value = 3.192871440749999
$ws.Range("S14").Value = 0.01631164764879812
$ws.Range("T14").Value = 0.01631164764879811

$ws.Range("G15").Value = 26.481475
$ws.Range("H15").Value = 79.444425
$ws.Range("I15").Value = 0.1943431960314218
$ws.Range("J15").Value = 0.1943431960314218
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.1260863333333333
$ws.Range("N15").Value = 0.378259
$ws.Range("O15").Value = 0.7899502755623498
$ws.Range("P15").Value = 0.7899502755623498
$ws.Range("Q15").Value = 3.338952084008333
$ws.Range("R15").Value = 30.050568756075
$ws.Range("S15").Value = 0.1535214612586894
$ws.Range("T15").Value = 0.1535214612586894

$ws.Range("G16").Value = 26.481475
$ws.Range("H16").Value = 79.444425
$ws.Range("I16").Value = 0.1943431960314218
$ws.Range("J16").Value = 0.1943431960314218
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.02013
$ws.Range("N16").Value = 0.06039
$ws.Range("O16").Value = 0.1261175468163621
$ws.Range("P16").Value = 0.1261175468163621
$ws.Range("Q16").Value = 0.5330720917499999
$ws.Range("R16").Value = 4.79764882575
$ws.Range("S16").Value = 0.02451008712393427
$ws.Range("T16").Value = 0.02451008712393427

